# Insert a new record at row 277 ("Feria Lagunitas de Puerto Montt" - Mango, fecha 2023-01-20),
# pushing the existing rows 277-305 down to 278-306 (the last existing row, old 305, ends up
# duplicated into the new row 306).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 277..305 down to 278..306, leaving row 277 blank (except inherited formatting).
$ws.Rows.Item(277).Insert()

# Populate the newly inserted row 277 with the new weekly record.
$ws.Range("A277").Value = 4
$ws.Range("B277").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C277").Value = "Los Lagos"
$ws.Range("D277").Value = 44946
$ws.Range("E277").Value = 10
$ws.Range("F277").Value = "Fruta"
$ws.Range("G277").Value = 100108
$ws.Range("H277").Value = "Tropicales y subtropicales"
$ws.Range("I277").Value = 100108002
$ws.Range("J277").Value = "Mango"
$ws.Range("K277").Value = "Sin especificar"
$ws.Range("L277").Value = "Primera"
$ws.Range("M277").Value = 200
$ws.Range("N277").Value = 7500
$ws.Range("O277").Value = 8000
$ws.Range("P277").Value = 7750
$ws.Range("Q277").Value = "`$/bandeja 4 kilos"
$ws.Range("R277").Value = "Brasil"
$ws.Range("S277").Value = 1938
$ws.Range("T277").Value = 4
